# Spike: exercise ActionSetting.target_slide & .hyperlink by building a
# first/previous/next/last slide-navigation deck around the existing slide.
#
# EMU -> point helper. PowerPoint COM stores Left/Top/Width/Height as
# (32-bit) points, and truncates when converting back to EMU on save, so a
# tiny epsilon is added to land on the exact target EMU value instead of
# 1 EMU short.
function EMU($emu) {
    return $emu / 914400 * 72 + 0.00004
}

$p = $ppt.ActivePresentation
$ppLayoutBlank = 12

# ---------------------------------------------------------------------
# Add the four navigation slides around the existing (only) slide, so the
# final order is: First, Previous, <original>, Next, Last.
# ---------------------------------------------------------------------

# 1) "First" slide, inserted before the existing slide.
$sFirst = $p.Slides.Add(1, $ppLayoutBlank)
$tbFirst = $sFirst.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tbFirst.TextFrame.WordWrap = 0
$tbFirst.TextFrame.AutoSize = 1
$tbFirst.TextFrame.TextRange.Text = "First"
$tbFirst.TextFrame2.VerticalAnchor = 3
$tbFirst.TextFrame.HorizontalAnchor = -1
$tbFirst.Left = EMU 4280157
$tbFirst.Top = EMU 3244334
$tbFirst.Width = EMU 583686
$tbFirst.Height = EMU 369332

# 2) "Previous" slide, inserted right after "First" (before the original).
$sPrev = $p.Slides.Add(2, $ppLayoutBlank)
$tbPrev = $sPrev.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tbPrev.TextFrame.WordWrap = 0
$tbPrev.TextFrame.AutoSize = 1
$tbPrev.TextFrame.TextRange.Text = "Previous"
$tbPrev.TextFrame2.VerticalAnchor = 3
$tbPrev.TextFrame.HorizontalAnchor = -1
$tbPrev.Left = EMU 4079397
$tbPrev.Top = EMU 3244334
$tbPrev.Width = EMU 985206
$tbPrev.Height = EMU 369332

# 3) "Next" slide, inserted right after the original slide (now at index 3).
$sNext = $p.Slides.Add(4, $ppLayoutBlank)
$tbNext = $sNext.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tbNext.TextFrame.WordWrap = 0
$tbNext.TextFrame.TextRange.Text = "Next"
$tbNext.Left = EMU 4260633
$tbNext.Top = EMU 3244334
$tbNext.Width = EMU 622735
$tbNext.Height = EMU 369332

# 4) "Last" slide, appended at the end.
$sLast = $p.Slides.Add(5, $ppLayoutBlank)
$tbLast = $sLast.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tbLast.TextFrame.WordWrap = 0
$tbLast.TextFrame.TextRange.Text = "Last"
$tbLast.Left = EMU 4293398
$tbLast.Top = EMU 3244334
$tbLast.Width = EMU 557204
$tbLast.Height = EMU 369332

Write-Output "slides now: $($p.Slides.Count)"
